# Generate Report for Archive
# - Status text moved from "Ready for handoff" to "In Translation" on all
#   three report sheets.
# - The Status/language columns are narrower now that the new status text
#   is shorter than the old one (column autofit shrank them).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status values that were "Ready for handoff".
$wsOverview.Range("E2:F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Shrink the columns that held the status text to match the new, shorter
# content (mirrors Excel auto-fitting these columns after the edit).
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
